# ConferenceParticipatedFaculty.xlsx — rework the Sheet1 header row:
#   A: Program Title            (unchanged)
#   B: School Name      -> Organizing Institute   (sample row: "school of xyz" -> "xyz")
#   C: Funded By         -> removed entirely (column deleted, shifting D/E left)
#   D->C: National / International -> Level       (data-validation list follows the column)
#   E->D: Year                                     (data-validation list follows the column)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "Funded By" column (C). Excel shifts D (National/International) and
# E (Year) left into C and D, carrying their values, styles and data validations along.
$ws.Columns.Item(3).Delete()

# New column B header + sample value.
$ws.Range("B1").Value = "Organizing Institute"
$ws.Range("B2").Value = "xyz"

# New column C header (was "National / International").
$ws.Range("C1").Value = "Level"

# Widen B and C to fit the renamed headers.
$ws.Columns.Item(2).ColumnWidth = 22.333333333333332
$ws.Columns.Item(3).ColumnWidth = 14.833333333333334

# Leave the same cell selected as in the edited workbook.
$ws.Range("E4").Select()
